$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same style (bold, bordered, centered)
# instead of Excel allocating a brand new style record.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-9
$values = @{
    2 = @(1, 2)
    3 = @(1, 5)
    4 = @(2, 6)
    5 = @(1, 4)
    6 = @(1, 2)
    7 = @(1, 4)
    8 = @(1, 3)
    9 = @(3, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
